# "Release the smoke test cases"
# - Swap ExecutionFlag values for row 7 (SMOKE) and row 8 (REGRESSION):
#   row7 A: No -> Yes
#   row8 A: Yes -> No
# - Move active selection from A9 to A8
# - Adjust window yWindow position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Yes"
$ws.Range("A8").Value = "No"

$ws.Range("A8").Select()

$excel.ActiveWindow.Top = 900
